# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 127 (pushing existing rows 127-135
# down to 128-136), then populate the new row with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(127).Insert()

$ws.Range("A127").Value2 = 8
$ws.Range("B127").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C127").Value2 = "Coquimbo"
$ws.Range("D127").Value2 = 44516
$ws.Range("E127").Value2 = 4
$ws.Range("F127").Value2 = 100112031
$ws.Range("G127").Value2 = "Poroto verde"
$ws.Range("H127").Value2 = "Magnum"
$ws.Range("I127").Value2 = "Primera"
$ws.Range("J127").Value2 = 400
$ws.Range("K127").Value2 = 42000
$ws.Range("L127").Value2 = 43000
$ws.Range("M127").Value2 = 42500
$ws.Range("N127").Value2 = "`$/caja 25 kilos"
$ws.Range("O127").Value2 = "Provincia de Limarí"
$ws.Range("P127").Value2 = 1700
$ws.Range("Q127").Value2 = 25
$ws.Range("R127").Value2 = "Hortaliza"
